# Refresh the crypto price/volume snapshot (scheduled GitHub Actions update).
# D-column price strings that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as text (preserving trailing
# zeros / locale-style "." thousands separators) instead of coercing them
# to a Number, matching the inline-string cells already in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.057.31"
$ws.Range("E2").Value = "  -1.83%  "

$ws.Range("D3").Value = "3.100.74"
$ws.Range("E3").Value = "  -3.38%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'212.92"
$ws.Range("E5").Value = "  -4.27%  "

$ws.Range("D6").Value = "'623.17"
$ws.Range("E6").Value = "  -2.96%  "

$ws.Range("D7").Value = "'0.377"
$ws.Range("E7").Value = "  -6.67%  "

$ws.Range("D8").Value = "'0.806"
$ws.Range("E8").Value = "  +13.82%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "3.102.82"
$ws.Range("E10").Value = "  -3.29%  "

$ws.Range("D11").Value = "'0.595"
$ws.Range("E11").Value = "  +2.49%  "

$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").Value = "'0.0000243"
$ws.Range("E13").Value = "  -6.91%  "

$ws.Range("D14").Value = "'5.30"
$ws.Range("E14").Value = "  -2.50%  "

$ws.Range("D15").Value = "88.529.90"
$ws.Range("E15").Value = "  -2.06%  "

$ws.Range("D16").Value = "'32.39"
$ws.Range("E16").Value = "  -4.00%  "

$ws.Range("D17").Value = "3.671.81"
$ws.Range("E17").Value = "  -3.25%  "

$ws.Range("D18").Value = "3.087.75"
$ws.Range("E18").Value = "  -3.56%  "

$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D20").Value = "'0.0000212"
$ws.Range("E20").Value = "  -7.11%  "

$ws.Range("D21").Value = "'13.48"
$ws.Range("E21").Value = "  -0.63%  "

$ws.Range("D22").Value = "'424.21"
$ws.Range("E22").Value = "  -4.00%  "

$ws.Range("D23").Value = "'8.29"
$ws.Range("E23").Value = "  -4.51%  "

$ws.Range("D24").Value = "'4.94"
$ws.Range("E24").Value = "  -2.94%  "

$ws.Range("D25").Value = "'5.66"
$ws.Range("E25").Value = "  +5.43%  "

$ws.Range("D26").Value = "'11.95"
$ws.Range("E26").Value = "  -0.31%  "

$ws.Range("D27").Value = "'82.60"
$ws.Range("E27").Value = "  +1.04%  "

$ws.Range("E28").Value = "  -4.25%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D32").Value = "'8.13"
$ws.Range("E32").Value = "  -4.28%  "

$ws.Range("D33").Value = "'512.82"
$ws.Range("E33").Value = "  -5.86%  "

$ws.Range("D34").Value = "'3.71"
$ws.Range("E34").Value = "  -13.04%  "

$ws.Range("D35").Value = "'6.82"
$ws.Range("E35").Value = "  -4.12%  "

$ws.Range("E36").Value = "  -3.51%  "

$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  -6.14%  "

$ws.Range("D38").Value = "'22.35"
$ws.Range("E38").Value = "  -1.12%  "

$ws.Range("D39").Value = "'22.28"
$ws.Range("E39").Value = "  -0.59%  "

$ws.Range("E40").Value = "  +1.62%  "

$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("E43").Value = "  -3.02%  "

$ws.Range("E44").Value = "  -6.28%  "

$ws.Range("D45").Value = "'146.32"
$ws.Range("E45").Value = "  -0.19%  "

$ws.Range("E46").Value = "  +4.58%  "

$ws.Range("D47").Value = "'0.0694"
$ws.Range("E47").Value = "  +13.48%  "

$ws.Range("D48").Value = "'43.47"
$ws.Range("E48").Value = "  -3.19%  "

$ws.Range("D49").Value = "'163.89"
$ws.Range("E49").Value = "  -5.93%  "

$ws.Range("D50").Value = "'1.23"
$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("D51").Value = "'0.708"
$ws.Range("E51").Value = "  -6.52%  "

# Rows 30-31: list order swapped - Binance-PegBSC-USD now ranks above Cronos
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.08"
$ws.Range("E30").Value = "  +7.72%  "

$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "'0.170"
$ws.Range("E31").Value = "  +6.64%  "
